$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New Water Year 2024 (and two more) ice-on records appended below the
# existing data block (rows 12-14), matching the style of column A (date).
$ws.Range("A12").Value = 45230
$ws.Range("B12").Value = 31
$ws.Range("C12").Value = 2024

$ws.Range("A13").Value = 45601
$ws.Range("B13").Value = 36
$ws.Range("C13").Value = 2025

$ws.Range("A14").Value = 45979
$ws.Range("B14").Value = 49
$ws.Range("C14").Value = 2026

# Carry the date number format down from the existing column A cells.
$ws.Range("A11").Copy()
$ws.Range("A12:A14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# Update the active selection to mirror where the author last clicked.
$ws.Range("B15").Select()
